$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G ("K") values recomputed (save_data regenerated to use K instead of Strike#)
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 2
    6  = 7
    7  = 6
    8  = 4
    9  = 2
    10 = 3
    11 = 1
    12 = 1
    13 = 1
    14 = 2
    15 = 2
    16 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
